# Finished Week 13 logging
# Adds two new RB players (R.Freeman, J.Samuels) with zeroed stats to the
# "RB" worksheet, and selects cell J7 afterward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RB")

# New player names to append below the existing roster (rows 2-4 are used).
$players = @("R.Freeman", "J.Samuels")

$row = 5
foreach ($name in $players) {
    $ws.Cells.Item($row, 1).Value = $name
    for ($col = 2; $col -le 10; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
    $row++
}

$ws.Range("J7").Select()
